$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Helper-ish constants (WdParagraphAlignment): 0=left 1=center 2=right 3=both(justify)
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# STEP 1: remove the four "Wenas" body paragraphs and the whole
# "El Problema de Investigación: El Habla Imaginada" chapter (title +
# body) together with its separating page-break paragraph, since none
# of that content survives in the edited document. Delete from the
# bottom up so earlier indices stay valid.
# ------------------------------------------------------------------
$toDelete = @(13, 12, 11, 10, 7, 4)
foreach ($idx in $toDelete) {
    $victim = $d.Paragraphs($idx)
    $victimRange = $d.Range($victim.Range.Start, $victim.Range.End)
    $victimRange.Delete()
}

# The document now contains exactly:
#   1 Marco Teórico
#   2 page break
#   3 "Una Introducción a la Neuroanatomía Humana"      (title -> chapter 1)
#   4 page break
#   5 "Electroencefalografía"                            (title -> chapter 4)
#   6 page break
#   7 "Inteligencia Artificial y Redes Neuronales"        (title -> chapter 5)

# ------------------------------------------------------------------
# STEP 2: Chapter 1 — reuse paragraph 3
# ------------------------------------------------------------------
$title1 = $d.Paragraphs(3).Range
[void]$title1.Find.Execute("Una Introducción a la Neuroanatomía Humana", $false, $false, $false, $false, $false, $true, 1, $false, "Soy un cerebro, Watson. El resto de mí es un mero apéndice", 2)
$d.Paragraphs(3).Format.Alignment = 1

# Make room for the epigraph + body paragraphs right after the title.
$d.Paragraphs(3).Range.InsertParagraphAfter()
$d.Paragraphs(4).Range.InsertParagraphAfter()

$epi1 = $d.Paragraphs(4)
$epi1.Range.InsertAfter("Arthur Conan Doyle, La piedra de Mazarino.")
$epi1 = $d.Paragraphs(4)
$epi1.Format.Alignment = 2
$epi1.Range.Font.Italic = $true
$epi1.Range.Font.ItalicBi = $true
$epi1.Range.Font.Size = 10
$epi1.Range.Font.SizeBi = 10

$body1 = $d.Paragraphs(5)
$body1.Range.InsertAfter("En este capítulo se hablará de una introducción a la neuroanatomía.")
$body1 = $d.Paragraphs(5)
$body1.Format.Alignment = 3
$body1.Range.Font.Size = 10
$body1.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# STEP 3: brand-new Chapter 2 — "Todo hombre puede ser..."
# Build it right after chapter 1's body (paragraph 5), before the
# existing page-break paragraph (now at index 6).
# ------------------------------------------------------------------
$anchor2 = $d.Paragraphs(5).Range
$anchor2.InsertParagraphAfter()
$anchor2.InsertParagraphAfter()
$anchor2.InsertParagraphAfter()
$anchor2.InsertParagraphAfter()

# paragraphs 6..9 are now blank placeholders: title, epigraph, body, page-break
$title2 = $d.Paragraphs(6)
$title2.Range.InsertAfter("Todo hombre puede ser, si se lo propone, escultor de su propio cerebro")
$title2 = $d.Paragraphs(6)
$title2.Format.Alignment = 3
$title2.Range.Font.Size = 14
$title2.Range.Font.SizeBi = 14

$epi2 = $d.Paragraphs(7)
$epi2.Range.InsertAfter("Santiago Ramón y Cajal.")
$epi2 = $d.Paragraphs(7)
$epi2.Format.Alignment = 2
$epi2.Range.Font.Italic = $true
$epi2.Range.Font.ItalicBi = $true
$epi2.Range.Font.Size = 10
$epi2.Range.Font.SizeBi = 10

$body2 = $d.Paragraphs(8)
$body2.Range.InsertAfter("En este capítulo se hablará centrándose ahora en los estudios de electroencefalografía.")
$body2 = $d.Paragraphs(8)
$body2.Format.Alignment = 3
$body2.Range.Font.Size = 10
$body2.Range.Font.SizeBi = 10

$pb2 = $d.Paragraphs(9)
$pb2.Range.Font.Size = 10
$pb2.Range.Font.SizeBi = 10
$r2 = $pb2.Range
$r2.Collapse(1)
$r2.InsertBreak(7)

# ------------------------------------------------------------------
# STEP 4: brand-new Chapter 3 — "Toda la tecnología tiende..."
# Inserted right after the page break we just created (paragraph 9),
# before the existing "Electroencefalografía" title (now at index 10).
# ------------------------------------------------------------------
$anchor3 = $d.Paragraphs(9).Range
$anchor3.InsertParagraphAfter()
$anchor3.InsertParagraphAfter()
$anchor3.InsertParagraphAfter()

$title3 = $d.Paragraphs(10)
$title3.Range.InsertAfter("Toda la tecnología tiende a crear un nuevo entorno humano")
$title3 = $d.Paragraphs(10)
$title3.Format.Alignment = 1
$title3.Range.Font.Size = 14
$title3.Range.Font.SizeBi = 14

$epi3 = $d.Paragraphs(11)
$epi3.Range.InsertAfter("Herbert Marshall Mcluhan")
$epi3 = $d.Paragraphs(11)
$epi3.Format.Alignment = 2
$epi3.Range.Font.Italic = $true
$epi3.Range.Font.ItalicBi = $true
$epi3.Range.Font.Size = 10
$epi3.Range.Font.SizeBi = 10

$body3 = $d.Paragraphs(12)
$body3.Range.InsertAfter("En este capítulo se hablará de la inteligencia artificial, las redes neuronales y todo el ámbito teórico necesario de esta área para entender correctamente los futuros capítulos.")
$body3 = $d.Paragraphs(12)
$body3.Format.Alignment = 3
$body3.Range.Font.Size = 10
$body3.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# STEP 5: Chapter 4 — reuse the "Electroencefalografía" title
# (now at paragraph 13, followed by its existing page-break paragraph)
# ------------------------------------------------------------------
$title4 = $d.Paragraphs(13).Range
[void]$title4.Find.Execute("Electroencefalografía", $false, $false, $false, $false, $false, $true, 1, $false, "El pensamiento humano puede literalmente, transformar el mundo físico", 2)
$d.Paragraphs(13).Format.Alignment = 3

$d.Paragraphs(13).Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.InsertParagraphAfter()

$epi4 = $d.Paragraphs(14)
$epi4.Range.InsertAfter("Dan Brown en su libro: El Símbolo Perdido.")
$epi4 = $d.Paragraphs(14)
$epi4.Format.Alignment = 2
$epi4.Range.Font.Italic = $true
$epi4.Range.Font.ItalicBi = $true
$epi4.Range.Font.Size = 10
$epi4.Range.Font.SizeBi = 10

$body4 = $d.Paragraphs(15)
$body4.Range.InsertAfter("Ha llegado el momento, pues tras todo el marco teórico previo, es el momento de abordar el problema de la presente investigación. Los temas previos funcionaron para poner sobre la mesa todo conocimiento necesario para enfocarse ahora en el habla imaginada.")
$body4 = $d.Paragraphs(15)
$body4.Format.Alignment = 3
$body4.Range.Font.Size = 10
$body4.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# STEP 6: Chapter 5 — reuse the "Inteligencia Artificial y Redes
# Neuronales" title (now shifted further down).
# ------------------------------------------------------------------
$title5 = $d.Paragraphs(17).Range
[void]$title5.Find.Execute("Inteligencia Artificial y Redes Neuronales", $false, $false, $false, $false, $false, $true, 1, $false, "Esas veces en las que estás muy cansado, no quieres exigirte pero igualmente lo haces… Ese es el sueño. No es el destino, es el trayecto", 2)
$d.Paragraphs(17).Format.Alignment = 1

$d.Paragraphs(17).Range.InsertParagraphAfter()
$d.Paragraphs(18).Range.InsertParagraphAfter()

$epi5 = $d.Paragraphs(18)
$epi5.Range.InsertAfter("Kobe Bryant")
$epi5 = $d.Paragraphs(18)
$epi5.Format.Alignment = 2
$epi5.Range.Font.Italic = $true
$epi5.Range.Font.ItalicBi = $true
$epi5.Range.Font.Size = 10
$epi5.Range.Font.SizeBi = 10

$body5 = $d.Paragraphs(19)
$body5.Range.InsertAfter("A continuación se presentará la experimentación propia realizada para buscar responder a la pregunta:")
$body5 = $d.Paragraphs(19)
$body5.Format.Alignment = 3
$body5.Range.Font.Size = 10
$body5.Range.Font.SizeBi = 10
